# covid-tracker-week.xlsx update
# - Increment existing id column (A2:A28) by 2 (rolling 14-day window shift)
# - Append a new record row (row 29) with id 101
# - Add a hyperlink for the new record's "source" cell (M29)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing ids (A2:A28) by +2 to make room for two newer records
# that rolled into the 14-day window.
for ($r = 2; $r -le 28; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value() + 2
}

# Append the new record as row 29.
$ws.Range("A29").Value = 101
$ws.Range("B29").Value = "patient_has_been_here"
$ws.Range("C29").Value = 44238
$ws.Range("C29").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D29").Value = "ตลาดท่าดินแดง คลองสาน"
$ws.Range("E29").Value = 13.73378448
$ws.Range("F29").Value = 100.5024518
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = " กทม.เปิดไทม์ไลน์ ผู้ป่วยโควิด  ประจำวันที่  19 กพ."
$ws.Hyperlinks.Add($ws.Range("M29"), "https://www.facebook.com/earthpongsakornk/posts/466640598079174")
$ws.Range("M29").Style = "Hyperlink"
$ws.Range("N29").Value = 0
